$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7799.8
$ws.Range("I62").Value = 8250
$ws.Range("J62").Value = 7499.6665
$ws.Range("K62").Value = 8250
$ws.Range("L62").Value = 7499.6665
$ws.Range("M62").Value = -7626
$ws.Range("N62").Value = -8747.666499999999

$ws.Range("H65").Value = 7799.8
$ws.Range("I65").Value = 8250
$ws.Range("J65").Value = 7499.6665
$ws.Range("K65").Value = 41250
$ws.Range("L65").Value = 37498.3325
$ws.Range("M65").Value = -38130
$ws.Range("N65").Value = -43738.3325

$ws.Range("H86").Value = 4324.5
$ws.Range("I86").Value = 1671
$ws.Range("K86").Value = 1671
$ws.Range("M86").Value = -548

$ws.Range("H89").Value = 4324.5
$ws.Range("I89").Value = 1671
$ws.Range("K89").Value = 8355
$ws.Range("M89").Value = -2739

$ws.Range("H92").Value = 818
$ws.Range("J92").Value = 798.75
$ws.Range("L92").Value = 798.75
$ws.Range("N92").Value = -3294.75

$ws.Range("H96").Value = 2327.8333
$ws.Range("I96").Value = 1193.4
$ws.Range("J96").Value = 8000
$ws.Range("K96").Value = 3580.2
$ws.Range("L96").Value = 24000
$ws.Range("M96").Value = -2207.2
$ws.Range("N96").Value = -26746

$ws.Range("H137").Value = 2074.0344
$ws.Range("J137").Value = 3191.6155
$ws.Range("L137").Value = 9574.8465
$ws.Range("N137").Value = -14674.8465

$ws.Range("H138").Value = 3619.9275
$ws.Range("I138").Value = 973.8461
$ws.Range("J138").Value = 4234.1963
$ws.Range("K138").Value = 2921.5383
$ws.Range("L138").Value = 12702.5889
$ws.Range("M138").Value = 2218.4617
$ws.Range("N138").Value = -22982.5889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 74.5
$ws.Range("I5").Value = 74.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 74.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 37.5
$ws.Range("N5").Value = $null

$ws.Range("H32").Value = 16485.291
$ws.Range("I32").Value = 7317.931
$ws.Range("J32").Value = 26710.424
$ws.Range("K32").Value = 7317.931
$ws.Range("L32").Value = 26710.424
$ws.Range("M32").Value = -7030.931
$ws.Range("N32").Value = -27284.424

$ws.Range("H74").Value = 4713.846
$ws.Range("I74").Value = 1293.4286
$ws.Range("K74").Value = 1293.4286
$ws.Range("M74").Value = -419.4286

$ws.Range("H77").Value = 4713.846
$ws.Range("I77").Value = 1293.4286
$ws.Range("K77").Value = 6467.143
$ws.Range("M77").Value = -2099.143

$ws.Range("H88").Value = 891
$ws.Range("I88").Value = 3199
$ws.Range("J88").Value = 314
$ws.Range("K88").Value = 3199
$ws.Range("L88").Value = 314
$ws.Range("M88").Value = -2793
$ws.Range("N88").Value = -1126

$ws.Range("H91").Value = 891
$ws.Range("I91").Value = 3199
$ws.Range("J91").Value = 314
$ws.Range("K91").Value = 3199
$ws.Range("L91").Value = 314
$ws.Range("M91").Value = -1795
$ws.Range("N91").Value = -3122

$ws.Range("H132").Value = 1003.73334
$ws.Range("I132").Value = 1003.73334
$ws.Range("K132").Value = 3011.20002
$ws.Range("M132").Value = -481.2000200000002

$ws.Range("H140").Value = 108999.5
$ws.Range("J140").Value = 108999.5
$ws.Range("L140").Value = 108999.5
$ws.Range("N140").Value = -119359.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 74.5
$ws.Range("I4").Value = 74.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 74.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 40.5
$ws.Range("N4").Value = $null

$ws.Range("H20").Value = 1720.6
$ws.Range("I20").Value = 1701.5714
$ws.Range("K20").Value = 1701.5714
$ws.Range("M20").Value = -1454.5714

$ws.Range("H86").Value = 4501.2856
$ws.Range("I86").Value = 4653.6
$ws.Range("J86").Value = 4416.6665
$ws.Range("K86").Value = 4653.6
$ws.Range("L86").Value = 4416.6665
$ws.Range("M86").Value = -3530.6
$ws.Range("N86").Value = -6662.6665

$ws.Range("H89").Value = 4501.2856
$ws.Range("I89").Value = 4653.6
$ws.Range("J89").Value = 4416.6665
$ws.Range("K89").Value = 23268
$ws.Range("L89").Value = 22083.3325
$ws.Range("M89").Value = -17652
$ws.Range("N89").Value = -33315.3325

$ws.Range("H107").Value = 1252
$ws.Range("I107").Value = 1265.8334
$ws.Range("K107").Value = 1265.8334
$ws.Range("M107").Value = 654.1666

$ws.Range("H134").Value = 3104.6365
$ws.Range("I134").Value = 2831.842
$ws.Range("K134").Value = 8495.526
$ws.Range("M134").Value = -5960.526

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3300
$ws.Range("I5").Value = 5000
$ws.Range("J5").Value = 1600
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 4800
$ws.Range("M5").Value = -14888
$ws.Range("N5").Value = -5024

$ws.Range("H129").Value = 2215.8667
$ws.Range("I129").Value = 1155.375
$ws.Range("K129").Value = 3466.125
$ws.Range("M129").Value = 1533.875

$ws.Range("H133").Value = 17514.5
$ws.Range("I133").Value = 15029
$ws.Range("K133").Value = 45087
$ws.Range("M133").Value = -40027

$ws.Range("H135").Value = 3300
$ws.Range("I135").Value = 5000
$ws.Range("J135").Value = 1600
$ws.Range("K135").Value = 45000
$ws.Range("L135").Value = 14400
$ws.Range("M135").Value = -42465
$ws.Range("N135").Value = -19470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 29000
$ws.Range("J53").Value = 29000
$ws.Range("L53").Value = 29000
$ws.Range("N53").Value = -30262

$ws.Range("H132").Value = 2499.3713
$ws.Range("I132").Value = 2144.6538
$ws.Range("K132").Value = 6433.9614
$ws.Range("M132").Value = -3903.9614

$ws.Range("H136").Value = 24560.834
$ws.Range("J136").Value = 24560.834
$ws.Range("L136").Value = 73682.50199999999
$ws.Range("N136").Value = -78782.50199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2165.6667
$ws.Range("I100").Value = 998
$ws.Range("K100").Value = 998
$ws.Range("M100").Value = -457

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").Value = $null

$ws.Range("H136").Value = 2223.7144
$ws.Range("I136").Value = 2223.7144
$ws.Range("K136").Value = 6671.1432
$ws.Range("M136").Value = -4121.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 113662.445
$ws.Range("I126").Value = 143709
$ws.Range("K126").Value = 431127
$ws.Range("M126").Value = -428657
